$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns from "physical quantity" terminology to "measurand" terminology
$ws.Range("F1").Value = "Measurand"
$ws.Range("I1").Value = "Measurand Level I"
$ws.Range("J1").Value = "Measurand Level II"

# Move the view/selection to F1, matching the saved sheet view state
$ws.Range("F1").Select()
